$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "centauro park shopping" row (row 2), shifting all rows below it up by one.
$ws.Rows.Item(2).Delete()
